# Updates the cryptos list (prices / 1h volume %) and swaps the
# Hedera / Aptos rows (38 and 39) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.422.97'
$ws.Range("E2").Value = '  +9.26%  '
$ws.Range("D3").Value = '1.681.87'
$ws.Range("E3").Value = '  +5.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9970'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3449'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.19'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.186'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07285'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9994'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("E13").Value = '  +4.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.140'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.767'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").Value = '1.677.33'
$ws.Range("E16").Value = '  +5.18%  '
$ws.Range("E17").Value = '  +3.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9967'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06730'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '81.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.39%  '
$ws.Range("E21").Value = '  +2.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.108'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '24.372.32'
$ws.Range("E24").Value = '  +9.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.432'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.687'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.358'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -11.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '1.860.98'
$ws.Range("E30").Value = '  +4.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '127.21'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.327'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.024'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9737'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.730'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08489'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.95%  '
$ws.Range("E37").Value = '  +5.44%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '12.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.32%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06504'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.361'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02344'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.262'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.92%  '
$ws.Range("E43").Value = '  +4.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6203'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9967'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.785'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5960'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '13.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.034'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07228'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.26%  '
